$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Ochai Agbaji (SG,SF / Toronto Raptors) -> Karl-Anthony Towns (PF,C / New York Knicks)
$ws.Range("A8").Value = "Karl-Anthony Towns"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "New York Knicks"

# Row 9: Karl-Anthony Towns (PF,C / New York Knicks) -> Julian Champagnie (SF,PF / San Antonio Spurs)
$ws.Range("A9").Value = "Julian Champagnie"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "San Antonio Spurs"

# Row 15: Trey Murphy III (SF,PF / New Orleans Pelicans) -> Austin Reaves (PG,SG / Los Angeles Lakers)
$ws.Range("A15").Value = "Austin Reaves"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Los Angeles Lakers"

# Row 16: Austin Reaves (PG,SG / Los Angeles Lakers) -> Mark Williams (C / Charlotte Hornets)
$ws.Range("A16").Value = "Mark Williams"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Charlotte Hornets"

# Row 17 stays: Franz Wagner (SF,PF / Orlando Magic) - unchanged

# Row 18: Mark Williams (C / Charlotte Hornets) -> Trey Murphy III (SF,PF / New Orleans Pelicans)
$ws.Range("A18").Value = "Trey Murphy III"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "New Orleans Pelicans"
